# Create and set up the environment for the Paging test suites
$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("login")

# Add the new "paging" sheet right after the existing "login" sheet.
$ws = $wb.Worksheets.Add($null, $login)
$ws.Name = "paging"

# Header row
$ws.Range("A1").Value = "*** Test Cases ***"
$ws.Range("B1").Value = "`${item}"
$ws.Range("C1").Value = "`${output}"
$ws.Range("D1").Value = "[Tags]"
$ws.Range("E1").Value = "[Documentation]"

# Row 2 - 9 items on a page
$ws.Range("A2").Value = "Check there are 9 items on 1 page"
$ws.Range("B2").Value = "Item on page 9"
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = "TC02"
$ws.Range("E2").Value = "Pagination: 9 items on a page"

# Row 3 - 18 items on a page
$ws.Range("A3").Value = "Check there are 18 items on 1 page"
$ws.Range("B3").Value = "Item on page 18"
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "TC02"
$ws.Range("E3").Value = "Pagination: 18 items on a page"

# Row 4 - 32 items on a page
$ws.Range("A4").Value = "Check there are 32 items on 1 page"
$ws.Range("B4").Value = "Item on page 32"
$ws.Range("C4").Value = 32
$ws.Range("D4").Value = "TC02"
$ws.Range("E4").Value = "Pagination: 32 items on a page"

# Row 5 - All items on a page
$ws.Range("A5").Value = "Check there are all items on 1 page"
$ws.Range("B5").Value = "All"
$ws.Range("C5").Value = "None"
$ws.Range("D5").Value = "TC02"
$ws.Range("E5").Value = "Pagination: All items on a page"

# Text-format the header row and the A/B/D/E columns (C stays General, like "login")
$ws.Range("A1:E1").NumberFormat = "@"
$ws.Range("A2:B5").NumberFormat = "@"
$ws.Range("D2:E5").NumberFormat = "@"

# Column widths for the new sheet (char widths, matching the "login" sheet look & feel)
$ws.Columns.Item(1).ColumnWidth = 29.830729166666668
$ws.Columns.Item(2).ColumnWidth = 16.276041666666668
$ws.Columns.Item(3).ColumnWidth = 9.608072916666666
$ws.Columns.Item(5).ColumnWidth = 50.385416666666664

[void]$ws.Range("B1").Select()

# Update the "login" sheet's view: drop the frozen/scrolled topLeftCell,
# keep it the selected tab, and move the active cell/selection to B4.
[void]$login.Select()
[void]$login.Range("B4").Select()
